# Update "Generate Report for Handback" timestamps across the workbook.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-07 09:42:32"

$wsZhCn.Range("H2").Value = "2016-09-07 09:42:26"
$wsZhCn.Range("K2").Value = "2016-09-07 09:42:46"

$wsDeDe.Range("H2").Value = "2016-09-07 09:42:32"
$wsDeDe.Range("K2").Value = "2016-09-07 09:42:55"
